$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 21743082
$ws.Range("I15").Value = 21743082
$ws.Range("K15").Value = 65229246
$ws.Range("M15").Value = -65229077
$ws.Range("H29").Value = 1302.9
$ws.Range("I29").Value = 99
$ws.Range("J29").Value = 1436.6666
$ws.Range("K29").Value = 297
$ws.Range("L29").Value = 4309.9998
$ws.Range("M29").Value = -16
$ws.Range("N29").Value = -4871.9998
$ws.Range("H111").Value = 1990.5714
$ws.Range("J111").Value = 1756.4
$ws.Range("L111").Value = 5269.200000000001
$ws.Range("N111").Value = -11403.2
$ws.Range("H132").Value = 2882.025
$ws.Range("I132").Value = 2888.9736
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 8666.9208
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -6136.9208
$ws.Range("N132").Value = -13310
$ws.Range("H138").Value = 210756.14
$ws.Range("I138").Value = 417444.4
$ws.Range("J138").Value = 4067.875
$ws.Range("K138").Value = 1252333.2
$ws.Range("L138").Value = 12203.625
$ws.Range("M138").Value = -1247193.2
$ws.Range("N138").Value = -22483.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8364.333000000001
$ws.Range("I37").Value = 3017
$ws.Range("J37").Value = 11038
$ws.Range("K37").Value = 3017
$ws.Range("L37").Value = 11038
$ws.Range("M37").Value = -2744
$ws.Range("N37").Value = -11584
$ws.Range("H61").Value = 1507
$ws.Range("I61").Value = 1364.8928
$ws.Range("J61").Value = 2833.3333
$ws.Range("K61").Value = 1364.8928
$ws.Range("L61").Value = 2833.3333
$ws.Range("M61").Value = -1152.8928
$ws.Range("N61").Value = -3257.3333
$ws.Range("H136").Value = 1507
$ws.Range("I136").Value = 1364.8928
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 4094.6784
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = -1544.6784
$ws.Range("N136").Value = -13599.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3076.4614
$ws.Range("I20").Value = 2813.8333
$ws.Range("J20").Value = 3301.5715
$ws.Range("K20").Value = 2813.8333
$ws.Range("L20").Value = 3301.5715
$ws.Range("M20").Value = -2566.8333
$ws.Range("N20").Value = -3795.5715
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H134").Value = 1339.8235
$ws.Range("I134").Value = 987.0857
$ws.Range("J134").Value = 2111.4375
$ws.Range("K134").Value = 2961.2571
$ws.Range("L134").Value = 6334.3125
$ws.Range("M134").Value = -426.2570999999998
$ws.Range("N134").Value = -11404.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 17859520
$ws.Range("I16").Value = 26318080
$ws.Range("J16").Value = 2558.111
$ws.Range("K16").Value = 26318080
$ws.Range("L16").Value = 2558.111
$ws.Range("M16").Value = -26317793
$ws.Range("N16").Value = -3132.111
$ws.Range("H22").Value = 492.66666
$ws.Range("I22").Value = 492.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 492.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -142.66666
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 4103484.8
$ws.Range("I31").Value = 6484842
$ws.Range("K31").Value = 6484842
$ws.Range("M31").Value = -6484547
$ws.Range("H34").Value = 4103484.8
$ws.Range("I34").Value = 6484842
$ws.Range("K34").Value = 6484842
$ws.Range("M34").Value = -6484640
$ws.Range("H113").Value = 17859520
$ws.Range("I113").Value = 26318080
$ws.Range("J113").Value = 2558.111
$ws.Range("K113").Value = 26318080
$ws.Range("L113").Value = 2558.111
$ws.Range("M113").Value = -26315910
$ws.Range("N113").Value = -6898.111
$ws.Range("H132").Value = 1295.1621
$ws.Range("I132").Value = 779.5925999999999
$ws.Range("J132").Value = 2687.2
$ws.Range("K132").Value = 2338.7778
$ws.Range("L132").Value = 8061.599999999999
$ws.Range("M132").Value = 191.2222000000002
$ws.Range("N132").Value = -13121.6
$ws.Range("H134").Value = 2601.1177
$ws.Range("I134").Value = 2783.3076
$ws.Range("J134").Value = 2009
$ws.Range("K134").Value = 8349.9228
$ws.Range("L134").Value = 6027
$ws.Range("M134").Value = -5814.9228
$ws.Range("N134").Value = -11097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1750.25
$ws.Range("J22").Value = 1750.25
$ws.Range("L22").Value = 5250.75
$ws.Range("N22").Value = -5588.75
$ws.Range("H27").Value = 1750.25
$ws.Range("J27").Value = 1750.25
$ws.Range("L27").Value = 5250.75
$ws.Range("N27").Value = -5454.75
$ws.Range("H58").Value = 66674770
$ws.Range("I58").Value = 3500
$ws.Range("J58").Value = 76931890
$ws.Range("K58").Value = 10500
$ws.Range("L58").Value = 230795670
$ws.Range("M58").Value = -10372
$ws.Range("N58").Value = -230795926
$ws.Range("H92").Value = 433.33334
$ws.Range("J92").Value = 450
$ws.Range("L92").Value = 1350
$ws.Range("N92").Value = -3846
$ws.Range("H110").Value = 10000
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H115").Value = 5101.3335
$ws.Range("I115").Value = 17209.334
$ws.Range("J115").Value = 3083.3333
$ws.Range("K115").Value = 51628.00199999999
$ws.Range("L115").Value = 9249.999899999999
$ws.Range("M115").Value = -50453.00199999999
$ws.Range("N115").Value = -11599.9999
$ws.Range("H121").Value = 1616.7587
$ws.Range("I121").Value = 561
$ws.Range("K121").Value = 1683
$ws.Range("M121").Value = -373
$ws.Range("H128").Value = 125000
$ws.Range("I128").Value = 125000
$ws.Range("K128").Value = 375000
$ws.Range("M128").Value = -370020
$ws.Range("H131").Value = 7093143
$ws.Range("J131").Value = 7247335.5
$ws.Range("L131").Value = 21742006.5
$ws.Range("N131").Value = -21752086.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1786.6471
$ws.Range("I126").Value = 1629.9166
$ws.Range("J126").Value = 2162.8
$ws.Range("K126").Value = 4889.7498
$ws.Range("L126").Value = 6488.400000000001
$ws.Range("M126").Value = -2419.7498
$ws.Range("N126").Value = -11428.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1826.9166
$ws.Range("I61").Value = 1577.35
$ws.Range("J61").Value = 3074.75
$ws.Range("K61").Value = 1577.35
$ws.Range("L61").Value = 3074.75
$ws.Range("M61").Value = -1375.35
$ws.Range("N61").Value = -3478.75
$ws.Range("H113").Value = 1826.9166
$ws.Range("I113").Value = 1577.35
$ws.Range("J113").Value = 3074.75
$ws.Range("K113").Value = 1577.35
$ws.Range("L113").Value = 3074.75
$ws.Range("M113").Value = 592.6500000000001
$ws.Range("N113").Value = -7414.75
$ws.Range("H132").Value = 10970618
$ws.Range("I132").Value = 19539562
$ws.Range("J132").Value = 2371.56
$ws.Range("K132").Value = 58618686
$ws.Range("L132").Value = 7114.68
$ws.Range("M132").Value = -58616156
$ws.Range("N132").Value = -12174.68
$ws.Range("H141").Value = 44440
$ws.Range("J141").Value = 44440
$ws.Range("L141").Value = 44440
$ws.Range("N141").Value = -54800
